$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original styling of the data range, then force Text format
# while assigning values, so numeric-looking strings (e.g. "1.001",
# "30.827.09") are stored as text instead of being auto-converted to numbers.
$dataRange = $ws.Range("B2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.827.09"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "1.928.15"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "241.37"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "0.4782"
$ws.Range("E7").Value = "  -1.84%  "

$ws.Range("D8").Value = "0.2880"
$ws.Range("E8").Value = "  -2.02%  "

$ws.Range("D9").Value = "0.06775"
$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("E10").Value = "  +2.35%  "

$ws.Range("D11").Value = "104.01"
$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("D12").Value = "0.07814"
$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").Value = "1.935.63"
$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("D15").Value = "0.6815"
$ws.Range("E15").Value = "  -2.49%  "

$ws.Range("D16").Value = "292.18"
$ws.Range("E16").Value = "  +7.91%  "

$ws.Range("D17").Value = "30.813.39"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000007580"
$ws.Range("E18").Value = "  -1.31%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.180.55"
$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("D21").Value = "12.87"
$ws.Range("E21").Value = "  -1.31%  "

$ws.Range("D22").Value = "5.513"
$ws.Range("E22").Value = "  -1.88%  "

$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("D24").Value = "6.377"
$ws.Range("E24").Value = "  -1.96%  "

$ws.Range("D25").Value = "9.534"
$ws.Range("E25").Value = "  -2.43%  "

$ws.Range("D26").Value = "168.52"
$ws.Range("E26").Value = "  +2.41%  "

$ws.Range("D27").Value = "19.77"
$ws.Range("E27").Value = "  +1.66%  "

$ws.Range("D28").Value = "2.110"
$ws.Range("E28").Value = "  -1.89%  "

$ws.Range("D29").Value = "1.391"
$ws.Range("E29").Value = "  +0.72%  "

$ws.Range("D30").Value = "0.1009"
$ws.Range("E30").Value = "  -2.27%  "

$ws.Range("D31").Value = "4.606"
$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("D32").Value = "1.531"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").Value = "4.343"
$ws.Range("E33").Value = "  -1.11%  "

$ws.Range("D34").Value = "0.04827"
$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("D35").Value = "0.7360"
$ws.Range("E35").Value = "  -1.93%  "

$ws.Range("E36").Value = "  -1.59%  "

$ws.Range("D37").Value = "2.726"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").Value = "0.01947"
$ws.Range("E38").Value = "  -2.33%  "

$ws.Range("D39").Value = "2.627"
$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("D40").Value = "6.443"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").Value = "75.33"
$ws.Range("E41").Value = "  -4.28%  "

$ws.Range("D42").Value = "2.024"
$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("D43").Value = "0.8666"
$ws.Range("E43").Value = "  -2.88%  "

$ws.Range("D44").Value = "0.4347"
$ws.Range("E44").Value = "  -1.48%  "

$ws.Range("D45").Value = "105.90"
$ws.Range("E45").Value = "  -1.79%  "

$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.003.18"
$ws.Range("E47").Value = "  +2.52%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.509"
$ws.Range("E48").Value = "  -4.28%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.116"
$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1208"
$ws.Range("E50").Value = "  -2.70%  "

$ws.Range("D51").Value = "35.01"
$ws.Range("E51").Value = "  -2.93%  "

# Restore original styling (so no stray style indices / number formats leak in)
$dataRange.Style = $origStyle
